# oxo testdata address, osprey EMEA UGC
# Adds 4 new "AccountDetails" rows (IL, WA, OK, CA) to the DataSet sheet,
# each mirroring the structure of row 2, with per-row address/city/state/zip,
# widens column B, and moves the view/selection down to the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")

$ppEmail = "LotusQA.GLD.PP.OXO.AutoTest2@gmail.com"
$prEmail = "LotusQA.GLD.PR.OXO.AutoTest2@gmail.com"
$pwd     = "Lotuswave@123"

# New account rows: AccountName, Street, City, State, Zip
$rows = @(
    @{ Row=39; Name="AccountDetails IL"; Street="3224 Sandy Ln";          City="Glenview"; State="Illinois";   Zip="60026" },
    @{ Row=40; Name="AccountDetails WA"; Street="419 Main St Unit 1";     City="Ione";      State="Washington"; Zip="99139" },
    @{ Row=41; Name="AccountDetails OK"; Street="508 N Grant St";         City="Cordell";   State="Oklahoma";   Zip="73632" },
    @{ Row=42; Name="AccountDetails CA"; Street="2309 Tulare St Unit 5";  City="Fresno";    State="California"; Zip="93707" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # --- Column A: account label ---
    $ws.Range("A${rowNum}").Value = $r.Name

    # --- Columns B/C/H/I: PP / PR emails (hyperlink-styled text) ---
    $ws.Range("B${rowNum}").Value = $ppEmail
    $ws.Range("C${rowNum}").Value = $prEmail
    $ws.Range("H${rowNum}").Value = $ppEmail
    $ws.Range("I${rowNum}").Value = $prEmail

    # --- Columns D/E: password (hyperlink-styled text) ---
    $ws.Range("D${rowNum}").Value = $pwd
    $ws.Range("E${rowNum}").Value = $pwd

    # --- Columns F/G: QA / TEST ---
    $ws.Range("F${rowNum}").Value = "QA"
    $ws.Range("G${rowNum}").Value = "TEST"

    # --- Address block: N Street, O City, P Country, Q/R State, S Zip ---
    $ws.Range("N${rowNum}").Value = $r.Street
    $ws.Range("O${rowNum}").Value = $r.City
    $ws.Range("P${rowNum}").Value = "United State"
    $ws.Range("Q${rowNum}").Value = $r.State
    $ws.Range("R${rowNum}").Value = $r.State
    $ws.Range("S${rowNum}").Value = "'" + $r.Zip

    # --- Column T: phone number (numeric) ---
    $ws.Range("T${rowNum}").Value = 9898989898

    # --- Column V: shipping method ---
    $ws.Range("V${rowNum}").Value = "Ground"

    # Apply the "Hyperlink" look to the cells that carry it in the template row
    # (also stamps J:M as blank-but-styled placeholder cells, matching row 2)
    $ws.Range("B${rowNum}:E${rowNum}").Style = "Hyperlink"
    $ws.Range("H${rowNum}:I${rowNum}").Style = "Hyperlink"
    $ws.Range("J${rowNum}:M${rowNum}").Style = "Hyperlink"

    # Real hyperlinks: E (password), C (PR email), I (PR email), B (PP email), H (PP email)
    $ws.Hyperlinks.Add($ws.Range("E${rowNum}"), "mailto:" + $pwd) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C${rowNum}"), "mailto:" + $prEmail) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I${rowNum}"), "mailto:" + $prEmail) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("B${rowNum}"), "mailto:" + $ppEmail) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("H${rowNum}"), "mailto:" + $ppEmail) | Out-Null

    # Re-apply the Hyperlink style after Hyperlinks.Add (which otherwise stamps its own xf)
    $ws.Range("B${rowNum}").Style = "Hyperlink"
    $ws.Range("C${rowNum}").Style = "Hyperlink"
    $ws.Range("E${rowNum}").Style = "Hyperlink"
    $ws.Range("H${rowNum}").Style = "Hyperlink"
    $ws.Range("I${rowNum}").Style = "Hyperlink"
}

# Widen column B (account email) and keep column C at its original width
$ws.Columns.Item(2).ColumnWidth = 25.28

# Move the view down to the newly-added rows
$ws.Activate()
$ws.Range("D45").Select()
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
